$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "97÷3=32, 1"
$t.Cell(1, 2).Range.Text = "70÷3=23, 1"
$t.Cell(1, 3).Range.Text = "15÷6=2, 3"
$t.Cell(1, 5).Range.Text = "68÷3=22, 2"
$t.Cell(5, 1).Range.Text = "10÷8=1, 2"
$t.Cell(5, 2).Range.Text = "60÷8=7, 4"
$t.Cell(5, 3).Range.Text = "24÷8=3, 0"
$t.Cell(5, 4).Range.Text = "96÷5=19, 1"
$t.Cell(5, 5).Range.Text = "54÷5=10, 4"
$t.Cell(9, 1).Range.Text = "28÷2=14, 0"
$t.Cell(9, 2).Range.Text = "53÷8=6, 5"
$t.Cell(9, 3).Range.Text = "83÷4=20, 3"
$t.Cell(9, 4).Range.Text = "19÷3=6, 1"
$t.Cell(9, 5).Range.Text = "46÷7=6, 4"
$t.Cell(13, 1).Range.Text = "29÷3=9, 2"
$t.Cell(13, 2).Range.Text = "45÷2=22, 1"
$t.Cell(13, 3).Range.Text = "55÷7=7, 6"
$t.Cell(13, 4).Range.Text = "64÷8=8, 0"
$t.Cell(13, 5).Range.Text = "55÷4=13, 3"
$t.Cell(17, 1).Range.Text = "33÷8=4, 1"
$t.Cell(17, 2).Range.Text = "64÷4=16, 0"
$t.Cell(17, 3).Range.Text = "80÷2=40, 0"
$t.Cell(17, 4).Range.Text = "71÷2=35, 1"
$t.Cell(17, 5).Range.Text = "23÷4=5, 3"
